# Commit: "Remove all non numeric characters from tel number"
# Applies to the "telephone" column (F) on the "customers 1" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers 1")
$ws.Activate()

$ws.Range("F2").Value = "+46731264413"
$ws.Range("F3").Value = "0046731212345"
$ws.Range("F5").Value = "004673 123-3332"

$ws.Range("G7").Select()
